# V = n_V_[T] fixed
# Bumps first_seed (col C) / last_seed (col D) for the two seed blocks on
# the "schedule" sheet:
#   rows 2-101   (T block 1): C 1 -> 21, D 10 -> 40
#   rows 102-201 (T block 2): C 11 -> 41, D 20 -> 60

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 101; $r++) {
    $ws.Cells.Item($r, 3).Value = 21
    $ws.Cells.Item($r, 4).Value = 40
}

for ($r = 102; $r -le 201; $r++) {
    $ws.Cells.Item($r, 3).Value = 41
    $ws.Cells.Item($r, 4).Value = 60
}

# Update the view state to match the saved workbook (scrolled pane +
# current selection).
$ws.Range("A43").Activate()
$ws.Range("I76").Select()
